# Generate Report for Handback
# Updates the localization-status workbook after a handback cycle:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" / "Latest Handback File"
#    columns populated (with a hyperlink on the target file), and the
#    "Latest Handback DateTime" stamped.
#  - The now-wider text in the Status / Target / Handback columns gets extra
#    column width so it isn't clipped.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceMd   = "3135611b-9e50-4dff-b862-c9ce08a97f02.md"
$sourceUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91b6467c582c6e61aa761ade7c5be8abf19c84fd/e2e/3135611b-9e50-4dff-b862-c9ce08a97f02.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet: status for both locales ----
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

# ---- zh-cn sheet ----
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C1").ColumnWidth = 29.17

$wsZhCn.Range("I2").Value = $sourceMd
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $sourceUrl, [Type]::Missing, $sourceMd, $sourceMd)
$wsZhCn.Range("I1").ColumnWidth = 39.17

$wsZhCn.Range("J2").Value = "3135611b-9e50-4dff-b862-c9ce08a97f02.c73e21a9a141ba993e077eee9cd57972ce52de00.zh-cn.xlf"
$wsZhCn.Range("J1").ColumnWidth = 39.17

$wsZhCn.Range("K2").Value = "2016-10-10 07:00:50"

# ---- de-de sheet ----
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C1").ColumnWidth = 29.17

$wsDeDe.Range("I2").Value = $sourceMd
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $sourceUrl, [Type]::Missing, $sourceMd, $sourceMd)
$wsDeDe.Range("I1").ColumnWidth = 39.17

$wsDeDe.Range("J2").Value = "3135611b-9e50-4dff-b862-c9ce08a97f02.c73e21a9a141ba993e077eee9cd57972ce52de00.de-de.xlf"
$wsDeDe.Range("J1").ColumnWidth = 39.17

$wsDeDe.Range("K2").Value = "2016-10-10 07:01:12"
